# multi browser implementation for chrome and firefox and accounting cash
$wb = $excel.ActiveWorkbook

$wsLoanInput = $wb.Worksheets.Item("NewLoanInput")
$wsSummary   = $wb.Worksheets.Item("Summary")
$wsSchedule  = $wb.Worksheets.Item("Repayment Schedule")

# --- NewLoanInput sheet: selection stays at G8 (tab no longer selected). ---
$wsLoanInput.Range("G8").Select() | Out-Null

# --- Summary sheet: move the selection from B16 to D4. ---
$wsSummary.Range("D4").Select() | Out-Null

# --- Repayment Schedule sheet: add a few blank-but-styled cells to row 2, ---
# --- matching the formatting already used by the neighboring B2 cell.    ---
$wsSchedule.Range("B2").Copy()
$wsSchedule.Range("A2").PasteSpecial(-4122)

$wsSchedule.Range("B2").Copy()
$wsSchedule.Range("E2").PasteSpecial(-4122)

$wsSchedule.Range("B2").Copy()
$wsSchedule.Range("N2").PasteSpecial(-4122)

$wsSchedule.Range("B2").Copy()
$wsSchedule.Range("O2").PasteSpecial(-4122)

# --- Update numeric values in row 2. ---
$wsSchedule.Range("H2").Value = 0
$wsSchedule.Range("J2").Value = $null
$wsSchedule.Range("K2").Value = 0
$wsSchedule.Range("L2").Value = 0

# --- Make Repayment Schedule the active tab/selection, replacing NewLoanInput. ---
$wsSchedule.Range("G11").Select() | Out-Null
$wsSchedule.Activate()
